# Update the solution values on the first results sheet (FTNC_Demand5).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FTNC_Demand5")

$ws.Range("B2").Value = 15.10426765927978
$ws.Range("C2").Value = 189.2671626970616
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 204.371430356341
